$wb = $excel.ActiveWorkbook

# --- Template sheet: clear the sample FR/NFR rows (now-empty template) ---
$wsTemplate = $wb.Worksheets.Item("Template")
$wsTemplate.Range("A3:C10").ClearContents()
$wsTemplate.Range("A3:M10").EntireRow.AutoFit()

# Update the title cell text (drop the "Current" date placeholder / mark it as the Template)
$wsTemplate.Range("A1").Value = "Template: Traceability Matrix" + [char]10 + "Project: " + [char]10 + "Date: "

# --- Example sheet: clear the second (FR/NFR) sample block, keep first block ---
$wsExample = $wb.Worksheets.Item("Example")
$wsExample.Range("A6:C13").ClearContents()
$wsExample.Range("A6:M13").EntireRow.AutoFit()

# --- Selection / active tab bookkeeping: Example becomes the active tab ---
$wsTemplate.Activate()
$wsTemplate.Range("A1:M1").Select()

$wsExample.Activate()
$wsExample.Range("C4").Select()
